# Slide 4 ("Content Placeholder 2") notes update:
#  - Paragraph "Track maker radius at jump base and peak." is removed.
#  - Paragraph "DangDxMax -> secLen" is removed.
#  - Paragraph "Suspension bottoming, ramp up stiffness ..." is removed.
#  - Paragraph "Momentum conservation across change in direction at jump
#    face transition." becomes the new first paragraph, prefixed with
#    "Check ".
#  - A new second paragraph is added describing the variable step solver.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item("Content Placeholder 2")
$tr = $shp.TextFrame.TextRange

# Replace the whole body with the (moved + edited) "Check Momentum ..." line,
# dropping the three paragraphs that used to precede it.
$tr.Text = "Check Momentum conservation across change in direction at jump face transition.  "

# Second paragraph: the new "Variable step solver?" note, built run-by-run so
# the inline code-style tokens (lwr_thr / upr_thr) stay separate runs, same
# as the rest of the deck's notes.
[void]$tr.InsertAfter("`rVariable step solver? As bike or wheel acceleration in either X or Y direction goes from [")
[void]$tr.InsertAfter("lwr_thr")
[void]$tr.InsertAfter(" -> ")
[void]$tr.InsertAfter("upr_thr")
[void]$tr.InsertAfter("], time steps goes [dt -> dt/n], where n is positive integer. Still only ")
[void]$tr.InsertAfter("record data every ")
[void]$tr.InsertAfter("dt, so run sub-loops at dt/n until a dt size step is complete, then record data")
